$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/20/2024  Through  5/26/2024"

# --- Weekly crime statistics table updates (rows 15-33) ---
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 6
$ws.Range("K15").Value = -33.333333333333
$ws.Range("L15").Value = -40
$ws.Range("M15").Value = -45.454545454545
$ws.Range("N15").Value = -57.142857142857
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 54.545454545454
$ws.Range("I16").Value = 81
$ws.Range("J16").Value = 79
$ws.Range("K16").Value = 2.53164556962
$ws.Range("L16").Value = 5.194805194805
$ws.Range("M16").Value = -13.829787234042
$ws.Range("N16").Value = -76.857142857142
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 29
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = 20.833333333333
$ws.Range("I17").Value = 140
$ws.Range("J17").Value = 125
$ws.Range("K17").Value = 12
$ws.Range("L17").Value = 26.126126126126
$ws.Range("M17").Value = 150
$ws.Range("N17").Value = -3.448275862068
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 49
$ws.Range("J18").Value = 39
$ws.Range("K18").Value = 25.641025641025
$ws.Range("L18").Value = -3.92156862745
$ws.Range("M18").Value = -53.77358490566
$ws.Range("N18").Value = -89.37093275488
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 37
$ws.Range("H19").Value = -32.727272727272
$ws.Range("I19").Value = 208
$ws.Range("J19").Value = 251
$ws.Range("K19").Value = -17.131474103585
$ws.Range("L19").Value = -27.017543859649
$ws.Range("M19").Value = 58.778625954198
$ws.Range("N19").Value = -8.370044052863
$ws.Range("C20").Value = 8
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 700
$ws.Range("F20").Value = 30
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 328.571428571429
$ws.Range("I20").Value = 130
$ws.Range("J20").Value = 82
$ws.Range("K20").Value = 58.536585365853
$ws.Range("L20").Value = 23.809523809523
$ws.Range("M20").Value = 11.111111111111
$ws.Range("N20").Value = -90.490124359912
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -10
$ws.Range("G21").Value = 103
$ws.Range("H21").Value = 18.446601941747
$ws.Range("I21").Value = 616
$ws.Range("J21").Value = 586
$ws.Range("K21").Value = 5.119453924914
$ws.Range("L21").Value = -3.75
$ws.Range("M21").Value = 19.37984496124
$ws.Range("N21").Value = -76.096235933255
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 9
$ws.Range("K22").Value = 12.5
$ws.Range("L22").Value = 80
$ws.Range("M22").Value = -10
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -40.54054054054
$ws.Range("F24").Value = 88
$ws.Range("G24").Value = 111
$ws.Range("H24").Value = -20.72072072072
$ws.Range("I24").Value = 486
$ws.Range("J24").Value = 518
$ws.Range("K24").Value = -6.177606177606
$ws.Range("L24").Value = -14.886164623467
$ws.Range("M24").Value = 89.84375
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 15
$ws.Range("E25").Value = -40
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 37
$ws.Range("H25").Value = 18.918918918918
$ws.Range("I25").Value = 238
$ws.Range("J25").Value = 211
$ws.Range("K25").Value = 12.796208530805
$ws.Range("L25").Value = -4.41767068273
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 9
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 49
$ws.Range("H26").Value = -6.122448979591
$ws.Range("I26").Value = 219
$ws.Range("J26").Value = 212
$ws.Range("K26").Value = 3.301886792452
$ws.Range("L26").Value = 23.033707865168
$ws.Range("M26").Value = 11.734693877551
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 10
$ws.Range("K27").Value = -23.076923076923
$ws.Range("L27").Value = -41.176470588235
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("I28").Value = 18
$ws.Range("J28").Value = 23
$ws.Range("K28").Value = -21.739130434782
$ws.Range("L28").Value = -28
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 2
$ws.Range("I29").Value = 2
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -60
$ws.Range("M29").Value = -33.333333333333
$ws.Range("N29").Value = -90.47619047619
$ws.Range("C30").Value = 1
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 2
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = -60
$ws.Range("M30").Value = -33.333333333333
$ws.Range("N30").Value = -90.47619047619
$ws.Range("F31").Value = 1
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = 50
$ws.Range("D33").Value = "0"
$ws.Range("E33").Value = "***.*"
